# GDPRPrj_TempiLavoro_LucaP.xlsx — "Inseriti form nel manuale"
# Adds a new timesheet row (34) for the GDPR "Documentazione" activity and a
# "Totale:" SUBTOTAL row (36) under the Tabella1 data, then scrolls/selects
# the view the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row -----------------------------------------------------
$ws.Range("E34").Value = "LucaP"
$ws.Range("F34").Value = "GDPR"
$ws.Range("G34").Value = "Documentazione"

# Copy the date formatting from the row above so H34 reuses the existing
# date style instead of minting a new number format, then set the value.
$ws.Range("H33").Copy()
$ws.Range("H34").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H34").Value = "3/22/2019"

$ws.Range("I34").Value = 70

# --- Totale row ---------------------------------------------------------
$ws.Range("H36").Value = "Totale:"
$ws.Range("I36").Formula = "=SUBTOTAL(109,I2:I35)"

# --- View state (scroll position + selection) ---------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
[void]$ws.Range("H36").Select()

Write-Host "done"
